# Auto-generated Excel COM-interop script to apply profit/price recalculation updates
# across the Excalibur_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 575
$ws.Range("I12").Value = 433.33334
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 433.33334
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = -263.33334
$ws.Range("N12").Value = -1340
$ws.Range("H15").Value = 2391.3333
$ws.Range("I15").Value = 2391.3333
$ws.Range("K15").Value = 7173.999899999999
$ws.Range("M15").Value = -7004.999899999999
$ws.Range("H18").Value = 2358.7646
$ws.Range("I18").Value = 925
$ws.Range("K18").Value = 925
$ws.Range("M18").Value = -641
$ws.Range("H40").Value = 7870.7144
$ws.Range("I40").Value = 8399
$ws.Range("J40").Value = 6550
$ws.Range("K40").Value = 8399
$ws.Range("L40").Value = 6550
$ws.Range("M40").Value = -8224
$ws.Range("N40").Value = -6900
$ws.Range("H113").Value = 4159.8
$ws.Range("H135").Value = 2757.25
$ws.Range("I135").Value = 3264.5
$ws.Range("K135").Value = 29380.5
$ws.Range("M135").Value = -26845.5
$ws.Range("H138").Value = 4276.0566
$ws.Range("J138").Value = 4869.7676
$ws.Range("L138").Value = 14609.3028
$ws.Range("N138").Value = -24889.3028

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1151.7333
$ws.Range("I32").Value = 1042.28
$ws.Range("J32").Value = 1699
$ws.Range("K32").Value = 1042.28
$ws.Range("L32").Value = 1699
$ws.Range("M32").Value = -755.28
$ws.Range("N32").Value = -2273
$ws.Range("H38").Value = 12633.333
$ws.Range("I38").Value = 12633.333
$ws.Range("K38").Value = 12633.333
$ws.Range("M38").Value = -12166.333
$ws.Range("H45").Value = 4164.067
$ws.Range("J45").Value = 2800.625
$ws.Range("L45").Value = 2800.625
$ws.Range("N45").Value = -3554.625
$ws.Range("H61").Value = 4238.8096
$ws.Range("I61").Value = 3501.75
$ws.Range("J61").Value = 6597.4
$ws.Range("K61").Value = 3501.75
$ws.Range("L61").Value = 6597.4
$ws.Range("M61").Value = -3289.75
$ws.Range("N61").Value = -7021.4
$ws.Range("H74").Value = 2802.84
$ws.Range("I74").Value = 1387.3529
$ws.Range("K74").Value = 1387.3529
$ws.Range("M74").Value = -513.3529000000001
$ws.Range("H77").Value = 2802.84
$ws.Range("I77").Value = 1387.3529
$ws.Range("K77").Value = 6936.7645
$ws.Range("M77").Value = -2568.7645
$ws.Range("H122").Value = 2533.7112
$ws.Range("I122").Value = 1835.1072
$ws.Range("K122").Value = 5505.321599999999
$ws.Range("M122").Value = -3055.321599999999
$ws.Range("H136").Value = 4238.8096
$ws.Range("I136").Value = 3501.75
$ws.Range("J136").Value = 6597.4
$ws.Range("K136").Value = 10505.25
$ws.Range("L136").Value = 19792.2
$ws.Range("M136").Value = -7955.25
$ws.Range("N136").Value = -24892.2
$ws.Range("H139").Value = 233059.58
$ws.Range("J139").Value = 233059.58
$ws.Range("L139").Value = 233059.58
$ws.Range("N139").Value = -243339.58

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 899.8333
$ws.Range("I22").Value = 899.8333
$ws.Range("K22").Value = 899.8333
$ws.Range("M22").Value = -726.8333
$ws.Range("H35").Value = 30000
$ws.Range("I35").Value = 30000
$ws.Range("K35").Value = 30000
$ws.Range("M35").Value = -29690
$ws.Range("H134").Value = 6988.5625
$ws.Range("I134").Value = 7175.3022
$ws.Range("K134").Value = 21525.9066
$ws.Range("M134").Value = -18990.9066

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 2460.2
$ws.Range("I11").Value = 2700.5
$ws.Range("J11").Value = 2300
$ws.Range("K11").Value = 2700.5
$ws.Range("L11").Value = 2300
$ws.Range("M11").Value = -2560.5
$ws.Range("N11").Value = -2580
$ws.Range("H22").Value = 757.6923
$ws.Range("I22").Value = 756.2
$ws.Range("K22").Value = 756.2
$ws.Range("M22").Value = -406.2
$ws.Range("H25").Value = 4424.75
$ws.Range("I25").Value = 2599.6667
$ws.Range("J25").Value = 9900
$ws.Range("K25").Value = 2599.6667
$ws.Range("L25").Value = 9900
$ws.Range("M25").Value = -2425.6667
$ws.Range("N25").Value = -10248
$ws.Range("H26").Value = 5142.857
$ws.Range("I26").Value = 5000
$ws.Range("J26").Value = 6250
$ws.Range("K26").Value = 5000
$ws.Range("L26").Value = 6250
$ws.Range("M26").Value = -4713
$ws.Range("N26").Value = -6824
$ws.Range("H58").Value = 3678.0908
$ws.Range("I58").Value = 1999
$ws.Range("K58").Value = 1999
$ws.Range("M58").Value = -1796
$ws.Range("H86").Value = 7982.3335
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 7982.3335
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 7982.3335
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -10228.3335
$ws.Range("H89").Value = 7982.3335
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 7982.3335
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 39911.6675
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -51143.6675
$ws.Range("H107").Value = 685.0454999999999
$ws.Range("I107").Value = 557.46155
$ws.Range("K107").Value = 557.46155
$ws.Range("M107").Value = 1362.53845
$ws.Range("H136").Value = 3678.0908
$ws.Range("I136").Value = 1999
$ws.Range("K136").Value = 5997
$ws.Range("M136").Value = -3447

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 557
$ws.Range("I19").Value = 499.5
$ws.Range("J19").Value = 580
$ws.Range("K19").Value = 1498.5
$ws.Range("L19").Value = 1740
$ws.Range("M19").Value = -1324.5
$ws.Range("N19").Value = -2088
$ws.Range("H114").Value = 3575.4092
$ws.Range("I114").Value = 1341.8334
$ws.Range("K114").Value = 4025.5002
$ws.Range("M114").Value = -771.5001999999999
$ws.Range("H118").Value = 1617
$ws.Range("I118").Value = 1234
$ws.Range("J118").Value = 2000
$ws.Range("K118").Value = 3702
$ws.Range("L118").Value = 6000
$ws.Range("M118").Value = -2459
$ws.Range("N118").Value = -8486
$ws.Range("H133").Value = 11333.333
$ws.Range("I133").Value = 9000
$ws.Range("K133").Value = 27000
$ws.Range("M133").Value = -21940

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 101769.6
$ws.Range("J13").Value = 2027.4286
$ws.Range("L13").Value = 2027.4286
$ws.Range("N13").Value = -2305.4286
$ws.Range("H97").Value = 3908.625
$ws.Range("I97").Value = 1799.5
$ws.Range("K97").Value = 1799.5
$ws.Range("M97").Value = -1303.5
$ws.Range("H102").Value = 2762.7856
$ws.Range("I102").Value = 1944.6875
$ws.Range("K102").Value = 1944.6875
$ws.Range("M102").Value = -322.6875
$ws.Range("H123").Value = 44666.668
$ws.Range("J123").Value = 44666.668
$ws.Range("L123").Value = 44666.668
$ws.Range("N123").Value = -49566.668
$ws.Range("H132").Value = 5599.4194
$ws.Range("I132").Value = 5081.826
$ws.Range("J132").Value = 7087.5
$ws.Range("K132").Value = 15245.478
$ws.Range("L132").Value = 21262.5
$ws.Range("M132").Value = -12715.478
$ws.Range("N132").Value = -26322.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3207.4546
$ws.Range("I7").Value = 3123.4375
$ws.Range("K7").Value = 3123.4375
$ws.Range("M7").Value = -3011.4375
$ws.Range("H40").Value = 1844.3334
$ws.Range("I40").Value = 2166.5
$ws.Range("J40").Value = 1200
$ws.Range("K40").Value = 2166.5
$ws.Range("L40").Value = 1200
$ws.Range("M40").Value = -2030.5
$ws.Range("N40").Value = -1472
$ws.Range("H43").Value = 20999.6
$ws.Range("J43").Value = 24999.5
$ws.Range("L43").Value = 24999.5
$ws.Range("N43").Value = -25385.5
$ws.Range("H46").Value = 3207.5833
$ws.Range("I46").Value = 1638.2
$ws.Range("J46").Value = 3620.5789
$ws.Range("K46").Value = 1638.2
$ws.Range("L46").Value = 3620.5789
$ws.Range("M46").Value = -1450.2
$ws.Range("N46").Value = -3996.5789
$ws.Range("H115").Value = 19000
$ws.Range("J115").Value = 19000
$ws.Range("L115").Value = 19000
$ws.Range("N115").Value = -21350
$ws.Range("H126").Value = 3207.4546
$ws.Range("I126").Value = 3123.4375
$ws.Range("K126").Value = 9370.3125
$ws.Range("M126").Value = -6900.3125
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 31651.53
$ws.Range("I136").Value = 3446.889
$ws.Range("K136").Value = 10340.667
$ws.Range("M136").Value = -7790.667000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 2579.6667
$ws.Range("I14").Value = 2108.6
$ws.Range("J14").Value = 4935
$ws.Range("K14").Value = 2108.6
$ws.Range("L14").Value = 4935
$ws.Range("M14").Value = -1940.6
$ws.Range("N14").Value = -5271
$ws.Range("H18").Value = 5000006
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H56").Value = 9142.5
$ws.Range("I56").Value = 9142.5
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 9142.5
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -8428.5
$ws.Range("N56").ClearContents()
$ws.Range("H100").Value = 1615.2142
$ws.Range("I100").Value = 1386.125
$ws.Range("K100").Value = 2772.25
$ws.Range("M100").Value = -2231.25
